$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RESUME
$ws2 = $wb.Worksheets.Item(2)   # RESULTATS

# ---------------------------------------------------------------------------
# 1) RESULTATS (sheet2): insert 4 new header columns (PASS / WARNING / FAIL /
#    ERROR) before the existing "Heure fin" column. Typed in an order that
#    reproduces the shared-string table ordering: PASS, FAIL, ERROR first,
#    then (later, once the RESUME sheet got its own "Warning" column) WARNING.
# ---------------------------------------------------------------------------
$ws2.Columns("C:F").Insert()
$ws2.Range("C1").Value = "PASS"
$ws2.Range("E1").Value = "FAIL"
$ws2.Range("F1").Value = "ERROR"
$ws2.Range("C1:F1").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 2) RESUME (sheet1): insert a "Warning" column in the results table (row 15
#    header / row 16 total) between "Passed" and "Failed".
# ---------------------------------------------------------------------------
$ws1.Columns("D").Insert()
$ws1.Range("D15").Value = "Warning"

# Back to RESULTATS: fill in the WARNING header now that the string exists.
$ws2.Range("D1").Value = "WARNING"

# Harmonise the header row (B15:F15) and the "Nombre de cas de tests" cell
# (A16) to the same (non coloured) bold style, and drop the old coloured
# per-status total cells (B16:E16) that used to sit next to it.
$ws1.Range("B15:F15").Font.Size = 12
$ws1.Range("B15:F15").Font.Bold = $true
$ws1.Range("A16").Font.Size = 12
$ws1.Range("B16:E16").Clear()

# Add a new "Nombre de STEP" row under "Nombre de cas de tests".
$ws1.Range("A17").Value = "Nombre de STEP"
$ws1.Range("A17").Font.Bold = $false
$ws1.Range("A17").Font.Size = 10
$ws1.Range("A17").HorizontalAlignment = -4152
$ws1.Range("A17").VerticalAlignment = -4108
$ws1.Rows("17").RowHeight = 24.6

# ---------------------------------------------------------------------------
# 3) Selections left by the editor when the workbook was saved.
# ---------------------------------------------------------------------------
$ws1.Range("C22").Select()
$ws2.Range("H16").Select()
